$d = $word.ActiveDocument
$bm = $d.Bookmarks.Item("_GoBack")
Write-Host "Bookmark exists: $($bm.Name)"
Write-Host "Bookmark start: $($bm.Start) end: $($bm.End)"
Write-Host "Bookmark range text: [$($bm.Range.Text)]"
